$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.781.86'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').Value = '3.697.65'
$ws.Range('E3').Value = '  +5.16%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '''419.39'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = '''130.49'
$ws.Range('E6').Value = '  -1.16%  '
$ws.Range('D7').Value = '3.688.64'
$ws.Range('E7').Value = '  +5.01%  '
$ws.Range('D8').Value = '''0.646'
$ws.Range('E8').Value = '  +1.45%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '''0.765'
$ws.Range('E10').Value = '  -2.81%  '
$ws.Range('D11').Value = '''0.183'
$ws.Range('E11').Value = '  +11.01%  '
$ws.Range('D12').Value = '''0.0000401'
$ws.Range('E12').Value = '  +50.68%  '
$ws.Range('D13').Value = '''43.12'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '''10.69'
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('D15').Value = '4.276.58'
$ws.Range('E15').Value = '  +4.82%  '
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '''20.63'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '3.692.98'
$ws.Range('E18').Value = '  +5.01%  '
$ws.Range('D19').Value = '''13.14'
$ws.Range('E19').Value = '  +4.56%  '
$ws.Range('D20').Value = '''1.13'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Value = '66.776.10'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').Value = '''443.43'
$ws.Range('E22').Value = '  -2.56%  '
$ws.Range('D23').Value = '''16.46'
$ws.Range('E23').Value = '  +22.75%  '
$ws.Range('D24').Value = '''90.28'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = '''3.16'
$ws.Range('E25').Value = '  -2.29%  '
$ws.Range('D26').Value = '''37.56'
$ws.Range('E26').Value = '  +9.70%  '
$ws.Range('D27').Value = '''10.24'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '''3.32'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('D29').Value = '''5.03'
$ws.Range('E29').Value = '  +4.74%  '
$ws.Range('D30').Value = '''0.127'
$ws.Range('E30').Value = '  +11.41%  '
$ws.Range('D31').Value = '''12.79'
$ws.Range('E31').Value = '  +1.79%  '
$ws.Range('D32').Value = '''2.78'
$ws.Range('E32').Value = '  -1.68%  '
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').Value = '''0.167'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').Value = '''41.57'
$ws.Range('E35').Value = '  +3.47%  '
$ws.Range('D36').Value = '''57.17'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('E38').Value = '  -3.01%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0740'
$ws.Range('E39').Value = '  +2.70%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '''0.153'
$ws.Range('E40').Value = '  +5.55%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').Value = '''3.08'
$ws.Range('E41').Value = '  +33.66%  '
$ws.Range('D42').Value = '''29.01'
$ws.Range('E42').Value = '  +31.04%  '
$ws.Range('D43').Value = '''0.998'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '''3.43'
$ws.Range('E44').Value = '  +2.43%  '
$ws.Range('D45').Value = '''149.13'
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('E46').Value = '  +4.77%  '
$ws.Range('D47').Value = '''2.89'
$ws.Range('E47').Value = '  -6.42%  '
$ws.Range('E48').Value = '  -5.09%  '
$ws.Range('E49').Value = '  -5.86%  '
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('E51').Value = '  +12.35%  '
